# New crime data collected - weekly CompStat 63rd Precinct update.
# Updates: report header (volume/number, week-covering dates) and the
# weekly/28-day/YTD crime-complaint statistics table (rows 15-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 30   Number  39" -> "...Number  40"
$ws.Range("A8").Value = "Volume 30   Number  40"

# --- Header: "Report Covering the Week  9/25/2023  Through  10/1/2023"
# --- becomes week of 10/2/2023 Through 10/8/2023
$ws.Range("C9").Value = "Report Covering the Week  10/2/2023  Through  10/8/2023"

# --- Row 15: Rape
$ws.Range("F15").Value = 1
$ws.Range("M15").Value = -33.333333333333

# --- Row 16: Robbery
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 90
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = -5.263157894736
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -47.976878612716
$ws.Range("N16").Value = -83.082706766917

# --- Row 17: Fel. Assault
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 131
$ws.Range("J17").Value = 136
$ws.Range("K17").Value = -3.676470588235
$ws.Range("L17").Value = 2.34375
$ws.Range("M17").Value = 29.702970297029
$ws.Range("N17").Value = -48.627450980392

# --- Row 18: Burglary
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 73
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = -19.780219780219
$ws.Range("L18").Value = -13.095238095238
$ws.Range("M18").Value = -66.046511627907
$ws.Range("N18").Value = -91.779279279279

# --- Row 19: Gr. Larceny
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 506
$ws.Range("J19").Value = 398
$ws.Range("K19").Value = 27.135678391959
$ws.Range("L19").Value = 80.714285714285
$ws.Range("M19").Value = 33.862433862433
$ws.Range("N19").Value = -1.171875

# --- Row 20: G.L.A.
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -7.692307692307
$ws.Range("I20").Value = 98
$ws.Range("J20").Value = 97
$ws.Range("K20").Value = 1.030927835051
$ws.Range("L20").Value = 60.655737704918
$ws.Range("M20").Value = -22.222222222222
$ws.Range("N20").Value = -95.416276894293

# --- Row 21: TOTAL
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -34.482758620689
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = -13.207547169811
$ws.Range("I21").Value = 911
$ws.Range("J21").Value = 824
$ws.Range("K21").Value = 10.558252427184
$ws.Range("L21").Value = 46.698872785829
$ws.Range("M21").Value = -9.712586719524
$ws.Range("N21").Value = -79.091117741565

# --- Row 23: Housing
$ws.Range("L23").Value = -8.695652173913
$ws.Range("M23").Value = -16

# --- Row 24: Petit Larceny
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 68.421052631578
$ws.Range("F24").Value = 124
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = 33.333333333333
$ws.Range("I24").Value = 990
$ws.Range("J24").Value = 879
$ws.Range("K24").Value = 12.627986348122
$ws.Range("L24").Value = 58.14696485623
$ws.Range("M24").Value = 28.571428571428

# --- Row 25: Misd. Assault
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 207
$ws.Range("J25").Value = 178
$ws.Range("K25").Value = 16.292134831460
$ws.Range("L25").Value = 8.376963350785
$ws.Range("M25").Value = -24.727272727272

# --- Row 26: UCR Rape*
$ws.Range("F26").Value = 1

# --- Row 27: Other Sex Crimes (C27 count drops to 0 -> rendered as text "0",
# matching the sheet's convention for zero counts, e.g. cells like C14/D14)
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = -8

$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# --- Row 28: Shooting Vic.
$ws.Range("M28").Value = -56.25

# --- Row 29: Shooting Inc.
$ws.Range("M29").Value = -53.846153846153
